$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1067
$ws.Range("I18").Value = 600.5
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 600.5
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -316.5
$ws.Range("N18").Value = -2568
$ws.Range("H62").Value = 2796
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2660
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2660
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -3908
$ws.Range("H65").Value = 2796
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2660
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 13300
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -19540
$ws.Range("H112").Value = 8153.4287
$ws.Range("J112").Value = 8617.272000000001
$ws.Range("L112").Value = 25851.816
$ws.Range("N112").Value = -28067.816
$ws.Range("H129").Value = 1400.3273
$ws.Range("I129").Value = 485.3125
$ws.Range("J129").Value = 1775.7179
$ws.Range("K129").Value = 1455.9375
$ws.Range("L129").Value = 5327.153700000001
$ws.Range("M129").Value = 3544.0625
$ws.Range("N129").Value = -15327.1537
$ws.Range("H138").Value = 2122.587
$ws.Range("I138").Value = 1352.7142
$ws.Range("J138").Value = 2769.28
$ws.Range("K138").Value = 4058.1426
$ws.Range("L138").Value = 8307.84
$ws.Range("M138").Value = 1081.8574
$ws.Range("N138").Value = -18587.84

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1727.2727
$ws.Range("J88").Value = 2333.3333
$ws.Range("L88").Value = 2333.3333
$ws.Range("N88").Value = -3145.3333
$ws.Range("H91").Value = 1727.2727
$ws.Range("J91").Value = 2333.3333
$ws.Range("L91").Value = 2333.3333
$ws.Range("N91").Value = -5141.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3728.5833
$ws.Range("I86").Value = 4740
$ws.Range("J86").Value = 2312.6
$ws.Range("K86").Value = 4740
$ws.Range("L86").Value = 2312.6
$ws.Range("M86").Value = -3617
$ws.Range("N86").Value = -4558.6
$ws.Range("H89").Value = 3728.5833
$ws.Range("I89").Value = 4740
$ws.Range("J89").Value = 2312.6
$ws.Range("K89").Value = 23700
$ws.Range("L89").Value = 11563
$ws.Range("M89").Value = -18084
$ws.Range("N89").Value = -22795

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 19333.867
$ws.Range("J4").Value = 19333.867
$ws.Range("L4").Value = 19333.867
$ws.Range("N4").Value = -19557.867
$ws.Range("H31").Value = 9215.714
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9215.714
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 9215.714
$ws.Range("N31").Value = -9805.714
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 9215.714
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9215.714
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9215.714
$ws.Range("N34").Value = -9619.714
$ws.Range("M34").ClearContents()
$ws.Range("H50").Value = 12999.25
$ws.Range("J50").Value = 12999.25
$ws.Range("L50").Value = 12999.25
$ws.Range("N50").Value = -14249.25
$ws.Range("H51").Value = 17399
$ws.Range("J51").Value = 17399
$ws.Range("L51").Value = 17399
$ws.Range("N51").Value = -18871
$ws.Range("H59").Value = 30399
$ws.Range("J59").Value = 30399
$ws.Range("L59").Value = 30399
$ws.Range("N59").Value = -32689
$ws.Range("H60").Value = 10448
$ws.Range("J60").Value = 10448
$ws.Range("L60").Value = 10448
$ws.Range("N60").Value = -11470
$ws.Range("H61").Value = 17399
$ws.Range("J61").Value = 17399
$ws.Range("L61").Value = 17399
$ws.Range("N61").Value = -18095
$ws.Range("H62").Value = 2952.7778
$ws.Range("I62").Value = 2675.7576
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 2675.7576
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -2051.7576
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 2952.7778
$ws.Range("I65").Value = 2675.7576
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 13378.788
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -10258.788
$ws.Range("N65").Value = -36240
$ws.Range("H68").Value = 23458.4
$ws.Range("J68").Value = 23458.4
$ws.Range("L68").Value = 23458.4
$ws.Range("N68").Value = -24956.4
$ws.Range("H71").Value = 23458.4
$ws.Range("J71").Value = 23458.4
$ws.Range("L71").Value = 70375.20000000001
$ws.Range("N71").Value = -77863.20000000001
$ws.Range("H74").Value = 29499.375
$ws.Range("J74").Value = 29499.375
$ws.Range("L74").Value = 29499.375
$ws.Range("N74").Value = -31247.375
$ws.Range("H77").Value = 29499.375
$ws.Range("J77").Value = 29499.375
$ws.Range("L77").Value = 88498.125
$ws.Range("N77").Value = -97234.125
$ws.Range("H122").Value = 1579.7179
$ws.Range("I122").Value = 1211.5294
$ws.Range("J122").Value = 1864.2273
$ws.Range("K122").Value = 3634.5882
$ws.Range("L122").Value = 5592.6819
$ws.Range("M122").Value = -1184.5882
$ws.Range("N122").Value = -10492.6819

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3092
$ws.Range("I131").Value = 383.5
$ws.Range("J131").Value = 4039.975
$ws.Range("K131").Value = 1150.5
$ws.Range("L131").Value = 12119.925
$ws.Range("M131").Value = 3889.5
$ws.Range("N131").Value = -22199.925
$ws.Range("H137").Value = 4908318.5
$ws.Range("I137").Value = 15163235
$ws.Range("J137").Value = 3793.4348
$ws.Range("K137").Value = 45489705
$ws.Range("L137").Value = 11380.3044
$ws.Range("M137").Value = -45484605
$ws.Range("N137").Value = -21580.3044
$ws.Range("H140").Value = 1938.5714
$ws.Range("I140").Value = 1392.1052
$ws.Range("K140").Value = 4176.3156
$ws.Range("M140").Value = 1003.6844

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3233.3333
$ws.Range("I122").Value = 2880
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8640
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6190
$ws.Range("N122").Value = -19900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10181.909
$ws.Range("I22").Value = 499.2
$ws.Range("J22").Value = 18250.834
$ws.Range("K22").Value = 499.2
$ws.Range("L22").Value = 18250.834
$ws.Range("M22").Value = -204.2
$ws.Range("N22").Value = -18840.834
$ws.Range("H27").Value = 10181.909
$ws.Range("I27").Value = 499.2
$ws.Range("J27").Value = 18250.834
$ws.Range("K27").Value = 499.2
$ws.Range("L27").Value = 18250.834
$ws.Range("M27").Value = -392.2
$ws.Range("N27").Value = -18464.834
$ws.Range("H132").Value = 2102.9524
$ws.Range("I132").Value = 1913.9595
$ws.Range("J132").Value = 3501.5
$ws.Range("K132").Value = 5741.8785
$ws.Range("L132").Value = 10504.5
$ws.Range("M132").Value = -3211.8785
$ws.Range("N132").Value = -15564.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 95000
$ws.Range("J95").Value = 95000
$ws.Range("L95").Value = 95000
$ws.Range("N95").Value = -100492
$ws.Range("H107").Value = 668.0270400000001
$ws.Range("I107").Value = 681.5517
$ws.Range("J107").Value = 619
$ws.Range("K107").Value = 2044.6551
$ws.Range("L107").Value = 1857
$ws.Range("M107").Value = -124.6550999999999
$ws.Range("N107").Value = -5697
$ws.Range("H136").Value = 1897.8636
$ws.Range("I136").Value = 1690.5625
$ws.Range("K136").Value = 5071.6875
$ws.Range("M136").Value = -2521.6875
